$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Year" column (D) stores its values as text in this workbook, not
# numbers. Force text formatting before assignment so Excel doesn't
# auto-convert the numeric-looking strings to numbers, then clear the
# formatting override again so the cell style matches the original
# (unstyled) data cells. (Row 5 is untouched by the edits below but is
# included here so the whole contiguous range round-trips identically.)
$yearCells = $ws.Range("D2:D11")
$yearCells.NumberFormat = "@"

# Row 2
$ws.Range("B2").Value = "Cyber scares and prophylactic policies: Crossnational evidence on the effect of cyberattacks on public support for surveillance"
$ws.Range("C2").Value = "Amelia C Arsenault, Sarah E Kreps, Keren LG Snider, Daphna Canetti"
$ws.Range("D2").Value = "2024"
$ws.Range("E2").Value = "10.1177/00223433241233960"

# Row 3
$ws.Range("B3").Value = "Warring from the virtual to the real: Assessing the public’s threshold for war over cyber security"
$ws.Range("C3").Value = "Sarah Kreps, Debak Das"
$ws.Range("D3").Value = "2017"
$ws.Range("E3").Value = "10.1177/2053168017715930"
$ws.Range("F3").Value = "Open Access"

# Row 4
$ws.Range("B4").Value = "The code not taken: China, the United States, and the future of cyber espionage"
$ws.Range("C4").Value = "Adam Segal"
$ws.Range("D4").Value = "2013"
$ws.Range("E4").Value = "10.1177/0096340213501344"

# Row 6
$ws.Range("B6").Value = "Framing cyber warfare: an analyst’s perspective"
$ws.Range("C6").Value = "Anthony Ween, Peter Dortmans, Nitin Thakur, Cayt Rowe"
$ws.Range("D6").Value = "2019"
$ws.Range("E6").Value = "10.1177/1548512917725620"
$ws.Range("F6").Value = "Restricted"

# Row 7
$ws.Range("B7").Value = "Mapping Global Cyberterror Networks: An Empirical Study of Al-Qaeda and ISIS Cyberterrorism Events"
$ws.Range("C7").Value = "Claire Seungeun Lee, Kyung-Shick Choi, Ryan Shandler, Chris Kayser"
$ws.Range("D7").Value = "2021"
$ws.Range("E7").Value = "10.1177/10439862211001606"

# Row 8
$ws.Range("B8").Value = "How the process of discovering cyberattacks biases our understanding of cybersecurity"
$ws.Range("C8").Value = "Harry Oppenheimer"
$ws.Range("D8").Value = "2024"
$ws.Range("E8").Value = "10.1177/00223433231217687"
$ws.Range("F8").Value = "Open Access"

# Row 9
$ws.Range("B9").Value = "Modeling Information Operations in a Tactical-level Stabilization Environment"
$ws.Range("C9").Value = "Helen Gaffney, Alasdair Vincent"
$ws.Range("D9").Value = "2011"
$ws.Range("E9").Value = "10.1177/1548512910388199"

# Row 10
$ws.Range("B10").Value = "Global versus Local Optimization in Redundancy Resolution of Robotic Manipulators"
$ws.Range("C10").Value = "Kazem Kazerounian, Zhaoyu Wang"
$ws.Range("D10").Value = "1988"
$ws.Range("E10").Value = "10.1177/027836498800700501"

# Row 11
$ws.Range("B11").Value = "Accountability and cyber conflict: examining institutional constraints on the use of cyber proxies"
$ws.Range("C11").Value = "William Akoto"
$ws.Range("D11").Value = "2022"
$ws.Range("E11").Value = "10.1177/07388942211051264"

# Restore the default (unstyled) appearance of the Year cells now that the
# text values are safely stored.
$yearCells.ClearFormats()
